$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update admin's password hash (row 2)
$ws.Range("B2").Value = '$2b$12$/or2YGXoZNVDkHa0iVcaGuZ7GK72fRUHAZCnpHcMF4IN2e8VwAXWC'

# Replace "diretor" row with "aluno" (row 3)
$ws.Range("A3").Value = "aluno"
$ws.Range("B3").Value = '$2b$12$L5aP6XHpk.FVIb0zVY5yR.8fttQPXKlcOy4U4u9qpwoTy32hsHW0W'

# Remove the old "professor" (row 4) and "aluno" (row 5) rows entirely
$ws.Rows("4:5").Delete()
